$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 8888
$ws.Range("J16").Value = 8888
$ws.Range("L16").Value = 8888
$ws.Range("N16").Value = -9348
$ws.Range("H74").Value = 4579.727
$ws.Range("I74").Value = 4100
$ws.Range("J74").Value = 5419.25
$ws.Range("K74").Value = 4100
$ws.Range("L74").Value = 5419.25
$ws.Range("M74").Value = -3164
$ws.Range("N74").Value = -7291.25
$ws.Range("H76").Value = 4905552.5
$ws.Range("I76").Value = 6947365.5
$ws.Range("J76").Value = 5200.3
$ws.Range("K76").Value = 6947365.5
$ws.Range("L76").Value = 5200.3
$ws.Range("M76").Value = -6947050.5
$ws.Range("N76").Value = -5830.3
$ws.Range("H77").Value = 4579.727
$ws.Range("I77").Value = 4100
$ws.Range("J77").Value = 5419.25
$ws.Range("K77").Value = 20500
$ws.Range("L77").Value = 27096.25
$ws.Range("M77").Value = -15820
$ws.Range("N77").Value = -36456.25
$ws.Range("H79").Value = 4905552.5
$ws.Range("I79").Value = 6947365.5
$ws.Range("J79").Value = 5200.3
$ws.Range("K79").Value = 6947365.5
$ws.Range("L79").Value = 5200.3
$ws.Range("M79").Value = -6946273.5
$ws.Range("N79").Value = -7384.3
$ws.Range("H137").Value = 1393.0714
$ws.Range("I137").Value = 1092.4517
$ws.Range("J137").Value = 2240.2727
$ws.Range("K137").Value = 3277.3551
$ws.Range("L137").Value = 6720.8181
$ws.Range("M137").Value = -727.3551000000002
$ws.Range("N137").Value = -11820.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 31890.889
$ws.Range("J23").Value = 9500
$ws.Range("L23").Value = 9500
$ws.Range("N23").Value = -10018
$ws.Range("H32").Value = 4716.3677
$ws.Range("I32").Value = 3400.1858
$ws.Range("K32").Value = 3400.1858
$ws.Range("M32").Value = -3113.1858
$ws.Range("H61").Value = 3044.7258
$ws.Range("I61").Value = 3154.509
$ws.Range("J61").Value = 2182.1428
$ws.Range("K61").Value = 3154.509
$ws.Range("L61").Value = 2182.1428
$ws.Range("M61").Value = -2942.509
$ws.Range("N61").Value = -2606.1428
$ws.Range("H136").Value = 3044.7258
$ws.Range("I136").Value = 3154.509
$ws.Range("J136").Value = 2182.1428
$ws.Range("K136").Value = 9463.527
$ws.Range("L136").Value = 6546.428400000001
$ws.Range("M136").Value = -6913.527
$ws.Range("N136").Value = -11646.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2400
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -1549
$ws.Range("N94").Value = -3902
$ws.Range("H105").Value = 31251960
$ws.Range("I105").Value = 45456160
$ws.Range("K105").Value = 45456160
$ws.Range("M105").Value = -45454413
$ws.Range("H134").Value = 3016.4558
$ws.Range("I134").Value = 3515.261
$ws.Range("J134").Value = 1973.5
$ws.Range("K134").Value = 10545.783
$ws.Range("L134").Value = 5920.5
$ws.Range("M134").Value = -8010.782999999999
$ws.Range("N134").Value = -10990.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 333335740
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 500003100
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 500003100
$ws.Range("M25").Value = -826
$ws.Range("N25").Value = -500003448
$ws.Range("H31").Value = 7816928
$ws.Range("I31").Value = 1623.8975
$ws.Range("J31").Value = 20008802
$ws.Range("K31").Value = 1623.8975
$ws.Range("L31").Value = 20008802
$ws.Range("M31").Value = -1328.8975
$ws.Range("N31").Value = -20009392
$ws.Range("H34").Value = 7816928
$ws.Range("I34").Value = 1623.8975
$ws.Range("J34").Value = 20008802
$ws.Range("K34").Value = 1623.8975
$ws.Range("L34").Value = 20008802
$ws.Range("M34").Value = -1421.8975
$ws.Range("N34").Value = -20009206
$ws.Range("H107").Value = 318
$ws.Range("I107").Value = 385.14285
$ws.Range("J107").Value = 291.8889
$ws.Range("K107").Value = 385.14285
$ws.Range("L107").Value = 291.8889
$ws.Range("M107").Value = 1534.85715
$ws.Range("N107").Value = -4131.8889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 401013.66
$ws.Range("I5").Value = 425
$ws.Range("J5").Value = 546682.25
$ws.Range("K5").Value = 1275
$ws.Range("L5").Value = 1640046.75
$ws.Range("M5").Value = -1163
$ws.Range("N5").Value = -1640270.75
$ws.Range("H11").Value = 157
$ws.Range("I11").Value = 81
$ws.Range("J11").Value = 195
$ws.Range("K11").Value = 243
$ws.Range("L11").Value = 585
$ws.Range("M11").Value = -103
$ws.Range("N11").Value = -865
$ws.Range("H19").Value = 2226.6667
$ws.Range("J19").Value = 2226.6667
$ws.Range("L19").Value = 6680.000100000001
$ws.Range("N19").Value = -7028.000100000001
$ws.Range("H22").Value = 2050.2
$ws.Range("J22").Value = 2362.75
$ws.Range("L22").Value = 7088.25
$ws.Range("N22").Value = -7426.25
$ws.Range("H27").Value = 2050.2
$ws.Range("J27").Value = 2362.75
$ws.Range("L27").Value = 7088.25
$ws.Range("N27").Value = -7292.25
$ws.Range("H122").Value = 4856.2
$ws.Range("I122").Value = 454.44446
$ws.Range("K122").Value = 4090.00014
$ws.Range("M122").Value = -1640.00014
$ws.Range("H135").Value = 401013.66
$ws.Range("I135").Value = 425
$ws.Range("J135").Value = 546682.25
$ws.Range("K135").Value = 3825
$ws.Range("L135").Value = 4920140.25
$ws.Range("M135").Value = -1290
$ws.Range("N135").Value = -4925210.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 11500000
$ws.Range("I20").Value = 11500000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 11500000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -11499755
$ws.Range("H51").Value = 50326
$ws.Range("J51").Value = 50326
$ws.Range("L51").Value = 50326
$ws.Range("N51").Value = -51344
$ws.Range("H80").Value = 2645.7693
$ws.Range("J80").Value = 2785.7144
$ws.Range("L80").Value = 2785.7144
$ws.Range("N80").Value = -4781.7144
$ws.Range("H83").Value = 2645.7693
$ws.Range("J83").Value = 2785.7144
$ws.Range("L83").Value = 13928.572
$ws.Range("N83").Value = -23912.572
$ws.Range("H122").Value = 8833910
$ws.Range("I122").Value = 7202980
$ws.Range("J122").Value = 12503502
$ws.Range("K122").Value = 21608940
$ws.Range("L122").Value = 37510506
$ws.Range("M122").Value = -21606490
$ws.Range("N122").Value = -37515406
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 7230
$ws.Range("J25").Value = 7230
$ws.Range("L25").Value = 7230
$ws.Range("N25").Value = -7816
$ws.Range("H132").Value = 1084.5714
$ws.Range("I132").Value = 735.63635
$ws.Range("J132").Value = 1956.909
$ws.Range("K132").Value = 2206.90905
$ws.Range("L132").Value = 5870.727000000001
$ws.Range("M132").Value = 323.0909499999998
$ws.Range("N132").Value = -10930.727
